$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correção das notas do fórum para matc65 em 2021.2:
# zera os indicadores diários (B:H), total_views (I) e nota_view (J)
# para todas as linhas de dados (linhas 2 a 50).
$ws.Range("B2:J50").Value = 0
